# "Final updation of the Code"
#
# Applies the geometry / text-formatting tweaks from the target diff:
#   - Slide 11 title "6. Result": turn on shrink-text-on-overflow
#     (<a:normAutofit/>) and bump the run to 44pt.
#   - Slide 13 picture: shorten its height (10781993 x 5430606 -> 10781993 x 5307514).
#   - Slide 8 picture (Content Placeholder 8): grow its height
#     (3193022 x 1754107 -> 3193022 x 4084068).
#   - Slide 9 picture (Content Placeholder 5): nudge its position by 1 EMU
#     and grow it (554477,1776193 / 6381344x4585696 -> 554476,1776192 / 6504435x4844415).
#
# (The diff also relocates the xmlns:p159 declaration from <mc:Choice> up to
#  <mc:AlternateContent> - and adds xmlns="" to <mc:Fallback> - on slides 12
#  and 13's transition markup. That's a namespace-prefix-only, semantically
#  inert rewrite of XML the PowerPoint object model doesn't expose: every
#  SlideShowTransition property setter either leaves that block's existing
#  serialization untouched or fully regenerates <p:transition>/its
#  mc:AlternateContent wrapper - dropping the morph fallback/adding
#  unrelated attributes - instead of just hoisting the namespace, so it is
#  intentionally left alone here rather than risk corrupting the transition.)


$p = $ppt.ActivePresentation

# --- Slide 11: "6. Result" title -------------------------------------------------
$s11 = $p.Slides.Item(11)
$title = $s11.Shapes.Item(1)
$title.TextFrame.AutoSize = 2          # ppAutoSizeTextToFitShape -> <a:normAutofit/>
$title.TextFrame.TextRange.Font.Size = 44

# --- Slide 13: resize the big screenshot -----------------------------------------
$s13 = $p.Slides.Item(13)
$pic13 = $s13.Shapes.Item(1)
$pic13.Height = 417.91448818897635

# --- Slide 8: resize the "Container" screenshot ----------------------------------
$s8 = $p.Slides.Item(8)
$pic8 = $s8.Shapes.Item(1)
$pic8.Height = 321.58015748031494

# --- Slide 9: reposition + resize the "Drop" screenshot --------------------------
$s9 = $p.Slides.Item(9)
$pic9 = $s9.Shapes.Item(1)
$pic9.Left = 43.65952875905512
$pic9.Top = 139.85764319527559
$pic9.Width = 512.1602478204725
$pic9.Height = 381.45
